# Auto-generated script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cellD = $ws.Range("D2")
$cellD.NumberFormat = "@"
$cellD.Value = "30.336.55"
$cellD.Style = "Normal"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3
$cellD = $ws.Range("D3")
$cellD.NumberFormat = "@"
$cellD.Value = "1.932.53"
$cellD.Style = "Normal"
$ws.Range("E3").Value = "  -0.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$cellD = $ws.Range("D5")
$cellD.NumberFormat = "@"
$cellD.Value = "0.7466"
$cellD.Style = "Normal"
$ws.Range("E5").Value = "  +2.89%  "

# Row 6
$cellD = $ws.Range("D6")
$cellD.NumberFormat = "@"
$cellD.Value = "248.53"
$cellD.Style = "Normal"
$ws.Range("E6").Value = "  -0.82%  "

# Row 7
$cellD = $ws.Range("D7")
$cellD.NumberFormat = "@"
$cellD.Value = "1.001"
$cellD.Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$cellD = $ws.Range("D8")
$cellD.NumberFormat = "@"
$cellD.Value = "28.29"
$cellD.Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "

# Row 9
$cellD = $ws.Range("D9")
$cellD.NumberFormat = "@"
$cellD.Value = "0.3205"
$cellD.Style = "Normal"
$ws.Range("E9").Value = "  -3.97%  "

# Row 10
$ws.Range("E10").Value = "  -2.14%  "

# Row 11
$cellD = $ws.Range("D11")
$cellD.NumberFormat = "@"
$cellD.Value = "0.7871"
$cellD.Style = "Normal"
$ws.Range("E11").Value = "  -3.25%  "

# Row 12
$cellD = $ws.Range("D12")
$cellD.NumberFormat = "@"
$cellD.Value = "0.08001"
$cellD.Style = "Normal"
$ws.Range("E12").Value = "  -1.20%  "

# Row 13
$cellD = $ws.Range("D13")
$cellD.NumberFormat = "@"
$cellD.Value = "1.936.43"
$cellD.Style = "Normal"
$ws.Range("E13").Value = "  -0.08%  "

# Row 14
$ws.Range("E14").Value = "  -2.00%  "

# Row 15
$cellD = $ws.Range("D15")
$cellD.NumberFormat = "@"
$cellD.Value = "94.53"
$cellD.Style = "Normal"
$ws.Range("E15").Value = "  -0.04%  "

# Row 16
$cellD = $ws.Range("D16")
$cellD.NumberFormat = "@"
$cellD.Value = "14.62"
$cellD.Style = "Normal"
$ws.Range("E16").Value = "  -2.67%  "

# Row 17
$cellD = $ws.Range("D17")
$cellD.NumberFormat = "@"
$cellD.Value = "30.341.00"
$cellD.Style = "Normal"
$ws.Range("E17").Value = "  -0.08%  "

# Row 18
$cellD = $ws.Range("D18")
$cellD.NumberFormat = "@"
$cellD.Value = "253.00"
$cellD.Style = "Normal"
$ws.Range("E18").Value = "  +1.04%  "

# Row 19
$cellD = $ws.Range("D19")
$cellD.NumberFormat = "@"
$cellD.Value = "0.000008030"
$cellD.Style = "Normal"
$ws.Range("E19").Value = "  -2.85%  "

# Row 20
$ws.Range("E20").Value = "  -1.82%  "

# Row 21
$cellD = $ws.Range("D21")
$cellD.NumberFormat = "@"
$cellD.Value = "2.187.40"
$cellD.Style = "Normal"
$ws.Range("E21").Value = "  -0.19%  "

# Row 22
$cellD = $ws.Range("D22")
$cellD.NumberFormat = "@"
$cellD.Value = "1.001"
$cellD.Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
$cellD = $ws.Range("D24")
$cellD.NumberFormat = "@"
$cellD.Value = "6.817"
$cellD.Style = "Normal"
$ws.Range("E24").Value = "  -2.28%  "

# Row 25
$cellD = $ws.Range("D25")
$cellD.NumberFormat = "@"
$cellD.Value = "9.574"
$cellD.Style = "Normal"
$ws.Range("E25").Value = "  -2.04%  "

# Row 26
$cellD = $ws.Range("D26")
$cellD.NumberFormat = "@"
$cellD.Value = "164.54"
$cellD.Style = "Normal"
$ws.Range("E26").Value = "  +0.77%  "

# Row 27
$cellD = $ws.Range("D27")
$cellD.NumberFormat = "@"
$cellD.Value = "2.334"
$cellD.Style = "Normal"
$ws.Range("E27").Value = "  -2.68%  "

# Row 28
$cellD = $ws.Range("D28")
$cellD.NumberFormat = "@"
$cellD.Value = "19.09"
$cellD.Style = "Normal"
$ws.Range("E28").Value = "  -1.09%  "

# Row 29
$cellD = $ws.Range("D29")
$cellD.NumberFormat = "@"
$cellD.Value = "0.1328"
$cellD.Style = "Normal"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30
$cellD = $ws.Range("D30")
$cellD.NumberFormat = "@"
$cellD.Value = "1.362"
$cellD.Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "

# Row 31
$ws.Range("E31").Value = "  -2.57%  "

# Row 32
$cellD = $ws.Range("D32")
$cellD.NumberFormat = "@"
$cellD.Value = "4.445"
$cellD.Style = "Normal"
$ws.Range("E32").Value = "  -0.17%  "

# Row 33
$cellD = $ws.Range("D33")
$cellD.NumberFormat = "@"
$cellD.Value = "4.146"
$cellD.Style = "Normal"
$ws.Range("E33").Value = "  -1.32%  "

# Row 34
$cellD = $ws.Range("D34")
$cellD.NumberFormat = "@"
$cellD.Value = "0.05140"
$cellD.Style = "Normal"
$ws.Range("E34").Value = "  -1.58%  "

# Row 35
$cellD = $ws.Range("D35")
$cellD.NumberFormat = "@"
$cellD.Value = "1.285"
$cellD.Style = "Normal"
$ws.Range("E35").Value = "  -0.87%  "

# Row 36
$cellD = $ws.Range("D36")
$cellD.NumberFormat = "@"
$cellD.Value = "0.7500"
$cellD.Style = "Normal"
$ws.Range("E36").Value = "  -0.12%  "

# Row 37
$cellD = $ws.Range("D37")
$cellD.NumberFormat = "@"
$cellD.Value = "2.773"
$cellD.Style = "Normal"
$ws.Range("E37").Value = "  +0.98%  "

# Row 38
$ws.Range("E38").Value = "  -0.81%  "

# Row 39
$cellD = $ws.Range("D39")
$cellD.NumberFormat = "@"
$cellD.Value = "2.802"
$cellD.Style = "Normal"
$ws.Range("E39").Value = "  -1.04%  "

# Row 40
$cellD = $ws.Range("D40")
$cellD.NumberFormat = "@"
$cellD.Value = "78.15"
$cellD.Style = "Normal"
$ws.Range("E40").Value = "  -3.53%  "

# Row 41
$ws.Range("E41").Value = "  -0.75%  "

# Row 42
$cellD = $ws.Range("D42")
$cellD.NumberFormat = "@"
$cellD.Value = "0.4507"
$cellD.Style = "Normal"
$ws.Range("E42").Value = "  -1.16%  "

# Row 43
$cellD = $ws.Range("D43")
$cellD.NumberFormat = "@"
$cellD.Value = "1.991"
$cellD.Style = "Normal"
$ws.Range("E43").Value = "  -2.80%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$cellD = $ws.Range("D44")
$cellD.NumberFormat = "@"
$cellD.Value = "1.001"
$cellD.Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cellD = $ws.Range("D45")
$cellD.NumberFormat = "@"
$cellD.Value = "0.8378"
$cellD.Style = "Normal"
$ws.Range("E45").Value = "  -1.21%  "

# Row 46
$cellD = $ws.Range("D46")
$cellD.NumberFormat = "@"
$cellD.Value = "102.76"
$cellD.Style = "Normal"
$ws.Range("E46").Value = "  +0.59%  "

# Row 47
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$cellD = $ws.Range("D47")
$cellD.NumberFormat = "@"
$cellD.Value = "7.573"
$cellD.Style = "Normal"
$ws.Range("E47").Value = "  +1.42%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cellD = $ws.Range("D48")
$cellD.NumberFormat = "@"
$cellD.Value = "9.820"
$cellD.Style = "Normal"
$ws.Range("E48").Value = "  +0.04%  "

# Row 49
$cellD = $ws.Range("D49")
$cellD.NumberFormat = "@"
$cellD.Value = "988.78"
$cellD.Style = "Normal"
$ws.Range("E49").Value = "  +12.53%  "

# Row 50
$cellD = $ws.Range("D50")
$cellD.NumberFormat = "@"
$cellD.Value = "37.48"
$cellD.Style = "Normal"
$ws.Range("E50").Value = "  +1.47%  "

# Row 51
$cellD = $ws.Range("D51")
$cellD.NumberFormat = "@"
$cellD.Value = "0.1213"
$cellD.Style = "Normal"
$ws.Range("E51").Value = "  +6.24%  "
